$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 1 for "Date and Time" (shifts all existing rows down by one).
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 13:17:50.432000 to 2024-03-11 15:29:55.824000"

# Insert a new row for "Cycle Count of battery" right before "Idling time percentage",
# which (after the previous insert) now sits at row 35.
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 75
